$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.284.39'
$ws.Range("E2").Value = '  +0.34%  '
$ws.Range("D3").Value = '1.590.68'
$ws.Range("E3").Value = '  +0.58%  '
$ws.Range("E4").Value = '  -0.22%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '212.82'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.58%  '
$ws.Range("E6").Value = '  +0.75%  '
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("E8").Value = '  +0.29%  '
$ws.Range("E9").Value = '  -0.22%  '
$ws.Range("E10").Value = '  -0.76%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0848'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.37%  '
$ws.Range("D12").Value = '1.814.50'
$ws.Range("E12").Value = '  +0.62%  '
$ws.Range("D13").Value = '1.610.04'
$ws.Range("E13").Value = '  +3.62%  '
$ws.Range("E14").Value = '  -0.06%  '
$ws.Range("E15").Value = '  +1.39%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '64.44'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = '26.282.39'
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("E18").Value = '  -0.98%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '7.47'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +2.44%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '213.45'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +2.82%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.29'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.89%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '9.01'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +1.41%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.15'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -2.14%  '
$ws.Range("E25").Value = '  +0.32%  '
$ws.Range("E27").Value = '  +0.87%  '
$ws.Range("E28").Value = '  -0.49%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '15.20'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  -0.92%  '
$ws.Range("E31").Value = '  +1.26%  '
$ws.Range("E32").Value = '  -0.19%  '
$ws.Range("E33").Value = '  +0.48%  '
$ws.Range("D34").Value = '1.336.84'
$ws.Range("E34").Value = '  +4.73%  '
$ws.Range("E35").Value = '  -0.90%  '
$ws.Range("E36").Value = '  -0.37%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.592'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -3.36%  '
$ws.Range("E38").Value = '  -0.14%  '
$ws.Range("E39").Value = '  +0.32%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '5.77'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +4.23%  '
$ws.Range("E41").Value = '  -0.22%  '
$ws.Range("E42").Value = '  -1.85%  '
$ws.Range("E43").Value = '  +0.37%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.762'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.39%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '61.87'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -0.66%  '
$ws.Range("D46").Value = '1.725.73'
$ws.Range("E46").Value = '  +0.54%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '86.69'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -2.66%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.50'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -3.33%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.0981'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -2.10%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0504'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.36%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.51%  '
